$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'boys small compression pants',
    'boys sport leggings',
    'boys sports tights leggings',
    'boys tights',
    'boys tights and leggings',
    'boys winter running pants',
    'boys winter running tights',
    'boys workout leggings',
    'brace for calf muscle',
    'brace for hamstring tear',
    'brace for runners knee',
    'braces essentials',
    'braces for spine',
    'braces pants',
    'bread rise basket',
    'breast compression',
    'breast reduction',
    'breath rocks',
    'breathable cloth tape',
    'breathable workout leggings for women',
    'bright blue tights',
    'bright capri leggings for women',
    'bright green tights for women',
    'bright key',
    'bright red leggings',
    'bright red tights',
    'bright tights for women',
    'bronchitis inhaler',
    'brooks greenlight tights',
    'brown athletic pants for women',
    'brown athletic tape',
    'brown capri tights',
    'brown compression',
    'brown compression pants',
    'brown compression tights',
    'brown design leggings',
    'brown leggings',
    'brown leggings for men',
    'brown leggins womens',
    'brown rise',
    'brown running pants',
    'brown thigh tights',
    'brown tights for women',
    'brown tights women',
    'brown tights womens',
    'brown winter dresses for women',
    'brown women',
    'brown womens tights',
    'brushed fleece lined leggings',
    'brushed fleece lined seamless leggings',
    'burlington women clothing',
    'burn xs',
    'butter blend underwear',
    'butter club shirts',
    'butter culture',
    'butter soft legging',
    'buttersoft leggings',
    'buttery soft black leggings',
    'buttery soft leggins',
    'buttocks support',
    'c3fit compression tights',
    'calf braces',
    'calf compression 2xu',
    'calf compression air',
    'calf compression band',
    'calf compression brace',
    'calf compression leggings',
    'calf compression pants',
    'calf compression running',
    'calf exercise',
    'calf injury support',
    'calf leggings for women',
    'calf muscle support',
    'calf muscle workout',
    'calf pro',
    'calf pro stretch',
    'calf running compression',
    'calf stretch',
    'calf support',
    'calf support band',
    'calf support brace 2 pack',
    'calf tightness',
    'calf workout',
    'calve compression',
    'calve machine',
    'calves compression',
    'calves ring',
    'calves workout',
    'calves workout equipment',
    'calzas deportivas para mujer',
    'calzas mujer deportiva nike',
    'came leggings for women',
    'camo nike',
    'camo pants women',
    'capri 3x',
    'capri activewear for women',
    'capri athletic leggings',
    'capri athletic leggings for women',
    'capri athletic pants',
    'capri bike shorts women'
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
